# Updates market-price-derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across several worksheets, as produced by the scheduled price-refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether / Ether
$ws.Range("H15").Value = 6007.9697
$ws.Range("I15").Value = 6007.9697
$ws.Range("K15").Value = 18023.9091
$ws.Range("M15").Value = -17854.9091

# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 8334979
$ws.Range("J17").Value = 9092614
$ws.Range("L17").Value = 27277842
$ws.Range("N17").Value = -27278178

# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 728.16
$ws.Range("I107").Value = 751
$ws.Range("K107").Value = 751
$ws.Range("M107").Value = 1169

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 4784.9287
$ws.Range("J116").Value = 5499.143
$ws.Range("L116").Value = 5499.143
$ws.Range("N116").Value = -12383.143

# Row 134: Binding Spells / Crocodileskin Index
$ws.Range("H134").Value = 92749.5
$ws.Range("J134").Value = 92749.5
$ws.Range("L134").Value = 92749.5
$ws.Range("N134").Value = -102889.5

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2298.7874
$ws.Range("I138").Value = 1968.7273
$ws.Range("J138").Value = 2399.639
$ws.Range("K138").Value = 5906.1819
$ws.Range("L138").Value = 7198.917
$ws.Range("M138").Value = -766.1818999999996
$ws.Range("N138").Value = -17478.917

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 2066.7917
$ws.Range("I2").Value = 956.1177
$ws.Range("J2").Value = 4764.143
$ws.Range("K2").Value = 956.1177
$ws.Range("L2").Value = 4764.143
$ws.Range("M2").Value = -843.1177
$ws.Range("N2").Value = -4990.143

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 3605.186
$ws.Range("I74").Value = 2647.3215
$ws.Range("J74").Value = 5393.2
$ws.Range("K74").Value = 2647.3215
$ws.Range("L74").Value = 5393.2
$ws.Range("M74").Value = -1773.3215
$ws.Range("N74").Value = -7141.2

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 3605.186
$ws.Range("I77").Value = 2647.3215
$ws.Range("J77").Value = 5393.2
$ws.Range("K77").Value = 13236.6075
$ws.Range("L77").Value = 26966
$ws.Range("M77").Value = -8868.6075
$ws.Range("N77").Value = -35702

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 2066.7917
$ws.Range("I116").Value = 956.1177
$ws.Range("J116").Value = 4764.143
$ws.Range("K116").Value = 956.1177
$ws.Range("L116").Value = 4764.143
$ws.Range("M116").Value = 1337.8823
$ws.Range("N116").Value = -9352.143

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 5548.07
$ws.Range("I132").Value = 4673.222
$ws.Range("K132").Value = 14019.666
$ws.Range("M132").Value = -11489.666

$ws = $wb.Worksheets.Item("BSM")
# Row 2: Proly Hatchet / Bronze Hatchet
$ws.Range("H2").Value = 64292.25
$ws.Range("J2").Value = 69056.336
$ws.Range("L2").Value = 69056.336
$ws.Range("N2").Value = -69282.336

# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 2066.7917
$ws.Range("I3").Value = 956.1177
$ws.Range("J3").Value = 4764.143
$ws.Range("K3").Value = 956.1177
$ws.Range("L3").Value = 4764.143
$ws.Range("M3").Value = -842.1177
$ws.Range("N3").Value = -4992.143

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2864.8542
$ws.Range("I134").Value = 2945.9565
$ws.Range("K134").Value = 8837.869499999999
$ws.Range("M134").Value = -6302.869499999999

$ws = $wb.Worksheets.Item("CRP")
# Row 21: Nightmare on My Street / Elm Cane
$ws.Range("H21").Value = 1250
$ws.Range("J21").Value = 1250
$ws.Range("L21").Value = 1250
$ws.Range("N21").Value = -1720

# Row 33: Tools for the Tools / Silver Battle Fork
$ws.Range("H33").Value = 348
$ws.Range("I33").Value = 348
$ws.Range("K33").Value = 348
$ws.Range("M33").Value = 31

# Row 35: Storm of Swords / Elm Macuahuitl
$ws.Range("H35").Value = 5061.8
$ws.Range("I35").Value = 3103
$ws.Range("J35").Value = 8000
$ws.Range("K35").Value = 3103
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = -2809
$ws.Range("N35").Value = -8588

# Row 51: Greenstone for Greenhorns / Jade Crook
$ws.Range("H51").Value = 21333.334
$ws.Range("I51").Value = 8000
$ws.Range("J51").Value = 28000
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 28000
$ws.Range("M51").Value = -7264
$ws.Range("N51").Value = -29472

# Row 61: Incant Now, Think Later / Jade Crook
$ws.Range("H61").Value = 21333.334
$ws.Range("I61").Value = 8000
$ws.Range("J61").Value = 28000
$ws.Range("K61").Value = 8000
$ws.Range("L61").Value = 28000
$ws.Range("M61").Value = -7652
$ws.Range("N61").Value = -28696

# Row 93: Reeling for Rods / Muudhorn Fishing Rod
$ws.Range("H93").Value = 2500
$ws.Range("I93").Value = 2500
$ws.Range("K93").Value = 2500
$ws.Range("M93").Value = -628

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 5683.8
$ws.Range("I99").Value = 5364.6665
$ws.Range("J99").Value = 6162.5
$ws.Range("K99").Value = 5364.6665
$ws.Range("L99").Value = 6162.5
$ws.Range("M99").Value = -3866.6665
$ws.Range("N99").Value = -9158.5

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 5683.8
$ws.Range("I126").Value = 5364.6665
$ws.Range("J126").Value = 6162.5
$ws.Range("K126").Value = 16093.9995
$ws.Range("L126").Value = 18487.5
$ws.Range("M126").Value = -13623.9995
$ws.Range("N126").Value = -23427.5

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 4852.543
$ws.Range("I134").Value = 3890.5173
$ws.Range("K134").Value = 11671.5519
$ws.Range("M134").Value = -9136.5519

$ws = $wb.Worksheets.Item("CUL")
# Row 58: Bread in the Clouds / La Noscean Toast
$ws.Range("H58").Value = 295
$ws.Range("I58").Value = 295
$ws.Range("K58").Value = 885
$ws.Range("M58").Value = -757

# Row 103: West Meats East / Nomad Meat Pie
$ws.Range("H103").Value = 1000
$ws.Range("I103").Value = 1000
$ws.Range("K103").Value = 3000
$ws.Range("M103").Value = -2121

# Row 106: Herky Jerky / Jerked Jhammel
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 8392.333000000001
$ws.Range("I2").Value = 67.90000000000001
$ws.Range("K2").Value = 67.90000000000001
$ws.Range("M2").Value = 45.09999999999999

# Row 39: One Man's Trash / Horn Ring
$ws.Range("H39").Value = 48000
$ws.Range("J39").Value = 50000
$ws.Range("L39").Value = 50000
$ws.Range("N39").Value = -51064

# Row 41: Renascence Man / Worm Fang Needle
$ws.Range("H41").Value = 43333
$ws.Range("I41").Value = 43333
$ws.Range("K41").Value = 43333
$ws.Range("M41").Value = -42978

# Row 69: High Above Me, She Sews Lovely / Mythrite Needle
$ws.Range("H69").Value = 33993
$ws.Range("J69").Value = 33993
$ws.Range("L69").Value = 33993
$ws.Range("N69").Value = -35491

# Row 72: Old-school Spooling (L) / Mythrite Needle
$ws.Range("H72").Value = 33993
$ws.Range("J72").Value = 33993
$ws.Range("L72").Value = 101979
$ws.Range("N72").Value = -109467

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 1773.6666
$ws.Range("I113").Value = 1808.6
$ws.Range("K113").Value = 1808.6
$ws.Range("M113").Value = 361.4000000000001

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 1931.3334
$ws.Range("I132").Value = 1772.0968
$ws.Range("J132").Value = 4399.5
$ws.Range("K132").Value = 5316.2904
$ws.Range("L132").Value = 13198.5
$ws.Range("M132").Value = -2786.2904
$ws.Range("N132").Value = -18258.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 2870.4285
$ws.Range("I40").Value = 2870.4285
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2870.4285
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2734.4285
$ws.Range("N40").Value = $null

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 201375.4
$ws.Range("J61").Value = 1984.5
$ws.Range("L61").Value = 1984.5
$ws.Range("N61").Value = -2388.5

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 201375.4
$ws.Range("J113").Value = 1984.5
$ws.Range("L113").Value = 1984.5
$ws.Range("N113").Value = -6324.5

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3947.8
$ws.Range("I122").Value = 3579.8333
$ws.Range("J122").Value = 4499.75
$ws.Range("K122").Value = 10739.4999
$ws.Range("L122").Value = 13499.25
$ws.Range("M122").Value = -8289.499899999999
$ws.Range("N122").Value = -18399.25

$ws = $wb.Worksheets.Item("WVR")
# Row 15: Workplace Safety / Cotton Scarf
$ws.Range("H15").Value = 9500
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

# Row 31: Whatchoo Talking About / Cotton Doublet Vest of Crafting
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = $null

# Row 74: Clothing the Naked Truth / Ramie Robe of Casting
$ws.Range("H74").Value = 19218.857
$ws.Range("I74").Value = 15969
$ws.Range("J74").Value = 19760.5
$ws.Range("K74").Value = 15969
$ws.Range("L74").Value = 19760.5
$ws.Range("M74").Value = -15033
$ws.Range("N74").Value = -21632.5

# Row 77: When in Robes (L) / Ramie Robe of Casting
$ws.Range("H77").Value = 19218.857
$ws.Range("I77").Value = 15969
$ws.Range("J77").Value = 19760.5
$ws.Range("K77").Value = 47907
$ws.Range("L77").Value = 59281.5
$ws.Range("M77").Value = -43227
$ws.Range("N77").Value = -68641.5

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 4538.5186
$ws.Range("I132").Value = 4138.227
$ws.Range("K132").Value = 12414.681
$ws.Range("M132").Value = -9884.681
